# Estadisticos Segundo Parcial 26 Mayo
# Update the gradebook summary rows (15, 17, 18, 19, 23) on the three
# sheets ("1er Parcial", "2o Parcial", "Final") with refreshed totals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "1er Parcial": only Promedio / Blancos / por_blancos move.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("1er Parcial")

$ws1.Range("I15").Value = 7.8
$ws1.Range("J15").Value = 0
$ws1.Range("K15").Value = 0

$ws1.Range("I17").Value = 8.6
$ws1.Range("J17").Value = 0
$ws1.Range("K17").Value = 0

$ws1.Range("J23").Value = 0
$ws1.Range("K23").Value = 0

# ---------------------------------------------------------------------
# Sheet "2o Parcial": Aprobados/Reprobados/percentages/Promedio/Blancos
# all refresh for rows 15, 17, 18, 19 and the summary row 23.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("2o Parcial")

$ws2.Range("E15").Value = 20
$ws2.Range("F15").Value = 3
$ws2.Range("G15").Value = 87
$ws2.Range("H15").Value = 13
$ws2.Range("I15").Value = 8.300000000000001
$ws2.Range("J15").Value = 0
$ws2.Range("K15").Value = 0

$ws2.Range("E17").Value = 37
$ws2.Range("F17").Value = 3
$ws2.Range("G17").Value = 92.5
$ws2.Range("H17").Value = 7.5
$ws2.Range("I17").Value = 8.699999999999999
$ws2.Range("J17").Value = 0
$ws2.Range("K17").Value = 0

$ws2.Range("E18").Value = 16
$ws2.Range("F18").Value = 1
$ws2.Range("G18").Value = 94.09999999999999
$ws2.Range("H18").Value = 5.9
$ws2.Range("I18").Value = 8.5
$ws2.Range("J18").Value = 0
$ws2.Range("K18").Value = 0

$ws2.Range("E19").Value = 16
$ws2.Range("F19").Value = 1
$ws2.Range("G19").Value = 94.09999999999999
$ws2.Range("H19").Value = 5.9
$ws2.Range("I19").Value = 8.5
$ws2.Range("J19").Value = 0
$ws2.Range("K19").Value = 0

$ws2.Range("E23").Value = 323
$ws2.Range("F23").Value = 34
$ws2.Range("G23").Value = 90.5
$ws2.Range("H23").Value = 9.5
$ws2.Range("I23").Value = 7.9
$ws2.Range("J23").Value = 0
$ws2.Range("K23").Value = 0

# ---------------------------------------------------------------------
# Sheet "Final": rows 15/17/23 get the same refreshed counts/percentages;
# rows 18/19 only move their Promedio column.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Final")

$ws3.Range("E15").Value = 20
$ws3.Range("F15").Value = 3
$ws3.Range("G15").Value = 87
$ws3.Range("H15").Value = 13
$ws3.Range("I15").Value = 8.1
$ws3.Range("J15").Value = 0
$ws3.Range("K15").Value = 0

$ws3.Range("E17").Value = 37
$ws3.Range("F17").Value = 3
$ws3.Range("G17").Value = 92.5
$ws3.Range("H17").Value = 7.5
$ws3.Range("I17").Value = 8.800000000000001
$ws3.Range("J17").Value = 0
$ws3.Range("K17").Value = 0

$ws3.Range("I18").Value = 8.1

$ws3.Range("I19").Value = 8.1

$ws3.Range("E23").Value = 323
$ws3.Range("F23").Value = 34
$ws3.Range("G23").Value = 90.5
$ws3.Range("H23").Value = 9.5
$ws3.Range("I23").Value = 8
$ws3.Range("J23").Value = 0
$ws3.Range("K23").Value = 0
